# Add files via upload
#
# The "data" sheet's column E (selection-reason color tag) is simplified from
# full color names to single-letter codes, and its header is shortened:
#   選點原因_Color -> Color
#   yellow         -> y
#   green          -> g
#   red            -> r
#   blue           -> b
#
# Sheet2 (the colour legend) keeps its values as-is.
#
# Finally the saved view/selection state is updated:
#   - "data" sheet: drop the scrolled topLeftCell, select H18
#   - "Sheet2"    : select B5 (and it is no longer the active tab)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item(2)

# Map the old full colour-name strings to their new single-letter codes.
$colorMap = @{
    "yellow" = "y"
    "green"  = "g"
    "red"    = "r"
    "blue"   = "b"
}

# Header rename.
$ws1.Range("E1").Value = "Color"

# Touch one representative cell of each colour first (in y, g, r, b order) so
# the new literals land in the shared-string table in that order, then sweep
# every data row to normalise the rest.
$ws1.Cells.Item(2, 5).Value   = "y"
$ws1.Cells.Item(61, 5).Value  = "g"
$ws1.Cells.Item(109, 5).Value = "r"
$ws1.Cells.Item(52, 5).Value  = "b"

$lastRow = $ws1.Cells.Item($ws1.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 141) { $lastRow = 141 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws1.Cells.Item($r, 5)
    $cur = [string]$cell.Text
    if ($colorMap.ContainsKey($cur)) {
        $cell.Value = $colorMap[$cur]
    }
}

# Restore view/selection state seen in the target workbook.
$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
$ws1.Range("H18").Select()
